$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing values (recomputed stats after 2024 river data refresh) ---
$ws.Range("G4").Value = 0.0214869285658454
$ws.Range("L4").Value = 0.00507
$ws.Range("G5").Value = 0.0214869285658454
$ws.Range("L5").Value = 0.00507
$ws.Range("G6").Value = 0.509260754375506
$ws.Range("L6").Value = 0.03926
$ws.Range("G7").Value = 0.509260754375506
$ws.Range("L7").Value = 0.03926
$ws.Range("G8").Value = 0.559342303620627
$ws.Range("L8").Value = 0.05056
$ws.Range("G9").Value = 0.559342303620627
$ws.Range("L9").Value = 0.05056
$ws.Range("G16").Value = 0.0217050965678054
$ws.Range("L16").Value = 0.00534
$ws.Range("G17").Value = 0.0217050965678054
$ws.Range("L17").Value = 0.00534
$ws.Range("G18").Value = 0.46971246652819
$ws.Range("L18").Value = 0.04837
$ws.Range("G19").Value = 0.46971246652819
$ws.Range("L19").Value = 0.04837
$ws.Range("G20").Value = 0.517612247394585
$ws.Range("L20").Value = 0.06172
$ws.Range("G21").Value = 0.517612247394585
$ws.Range("L21").Value = 0.06172
$ws.Range("G28").Value = 1147.48903097379
$ws.Range("H28").Value = 12981.0063938204
$ws.Range("G29").Value = 1147.48903097379
$ws.Range("H29").Value = 12981.0063938204
$ws.Range("G30").Value = 1147.48903097379
$ws.Range("H30").Value = 12981.0063938204
$ws.Range("G31").Value = 1147.48903097379
$ws.Range("H31").Value = 12981.0063938204
$ws.Range("G32").Value = 0.035285702905245
$ws.Range("G33").Value = 0.035285702905245
$ws.Range("G34").Value = 0.513611716624731
$ws.Range("G35").Value = 0.513611716624731
$ws.Range("G36").Value = 0.578233401043477
$ws.Range("G37").Value = 0.578233401043477

# --- Append new attribute rows 42-57 for period 2019 - 2023 ---
# Row 42
$ws.Cells.Item(42,1).Value = "Kaitoke at Vector Gas Line"
$ws.Cells.Item(42,2).Value = "DRP (95th Percentile)"
$ws.Cells.Item(42,3).Value = "D"
$ws.Cells.Item(42,4).Value = "2019 - 2023"
$ws.Cells.Item(42,5).Value = "RepSite"
$ws.Cells.Item(42,6).Value = 0.0375
$ws.Cells.Item(42,7).Value = 0.0460555555555556
$ws.Cells.Item(42,8).Value = 0.202
$ws.Cells.Item(42,9).Value = 0.119
$ws.Cells.Item(42,12).Value = 0.048
$ws.Cells.Item(42,13).Value = 0.06032
$ws.Cells.Item(42,14).Value = 0.06772
$ws.Cells.Item(42,15).Value = 1773468
$ws.Cells.Item(42,16).Value = 5573594
$ws.Cells.Item(42,17).Value = "Whanganui District"
$ws.Cells.Item(42,18).Value = "Whanganui"
$ws.Cells.Item(42,19).Value = "Kaitoke Lakes"
$ws.Cells.Item(42,20).Value = "West_4"
$ws.Cells.Item(42,21).Value = "mg/L"

# Row 43
$ws.Cells.Item(43,1).Value = "Kaitoke at Vector Gas Line"
$ws.Cells.Item(43,2).Value = "DRP (Median)"
$ws.Cells.Item(43,3).Value = "D"
$ws.Cells.Item(43,4).Value = "2019 - 2023"
$ws.Cells.Item(43,5).Value = "RepSite"
$ws.Cells.Item(43,6).Value = 0.0375
$ws.Cells.Item(43,7).Value = 0.0460555555555556
$ws.Cells.Item(43,8).Value = 0.202
$ws.Cells.Item(43,9).Value = 0.119
$ws.Cells.Item(43,12).Value = 0.048
$ws.Cells.Item(43,13).Value = 0.06032
$ws.Cells.Item(43,14).Value = 0.06772
$ws.Cells.Item(43,15).Value = 1773468
$ws.Cells.Item(43,16).Value = 5573594
$ws.Cells.Item(43,17).Value = "Whanganui District"
$ws.Cells.Item(43,18).Value = "Whanganui"
$ws.Cells.Item(43,19).Value = "Kaitoke Lakes"
$ws.Cells.Item(43,20).Value = "West_4"
$ws.Cells.Item(43,21).Value = "mg/L"

# Row 44
$ws.Cells.Item(44,1).Value = "Kaitoke at Vector Gas Line"
$ws.Cells.Item(44,2).Value = "E coli (>260)"
$ws.Cells.Item(44,3).Value = "D"
$ws.Cells.Item(44,4).Value = "2019 - 2023"
$ws.Cells.Item(44,5).Value = "RepSite"
$ws.Cells.Item(44,6).Value = 250
$ws.Cells.Item(44,7).Value = 931.722340626303
$ws.Cells.Item(44,8).Value = 12981.0063938204
$ws.Cells.Item(44,9).Value = 5860
$ws.Cells.Item(44,10).Value = 27.7777777777778
$ws.Cells.Item(44,11).Value = 48.1481481481481
$ws.Cells.Item(44,12).Value = 222
$ws.Cells.Item(44,13).Value = 740.16
$ws.Cells.Item(44,14).Value = 2696
$ws.Cells.Item(44,15).Value = 1773468
$ws.Cells.Item(44,16).Value = 5573594
$ws.Cells.Item(44,17).Value = "Whanganui District"
$ws.Cells.Item(44,18).Value = "Whanganui"
$ws.Cells.Item(44,19).Value = "Kaitoke Lakes"
$ws.Cells.Item(44,20).Value = "West_4"
$ws.Cells.Item(44,21).Value = "% exceedances over 260/100 mL"

# Row 45
$ws.Cells.Item(45,1).Value = "Kaitoke at Vector Gas Line"
$ws.Cells.Item(45,2).Value = "E coli (>540)"
$ws.Cells.Item(45,3).Value = "D"
$ws.Cells.Item(45,4).Value = "2019 - 2023"
$ws.Cells.Item(45,5).Value = "RepSite"
$ws.Cells.Item(45,6).Value = 250
$ws.Cells.Item(45,7).Value = 931.722340626303
$ws.Cells.Item(45,8).Value = 12981.0063938204
$ws.Cells.Item(45,9).Value = 5860
$ws.Cells.Item(45,10).Value = 27.7777777777778
$ws.Cells.Item(45,11).Value = 48.1481481481481
$ws.Cells.Item(45,12).Value = 222
$ws.Cells.Item(45,13).Value = 740.16
$ws.Cells.Item(45,14).Value = 2696
$ws.Cells.Item(45,15).Value = 1773468
$ws.Cells.Item(45,16).Value = 5573594
$ws.Cells.Item(45,17).Value = "Whanganui District"
$ws.Cells.Item(45,18).Value = "Whanganui"
$ws.Cells.Item(45,19).Value = "Kaitoke Lakes"
$ws.Cells.Item(45,20).Value = "West_4"
$ws.Cells.Item(45,21).Value = "% exceedances over 540/100 mL"

# Row 46
$ws.Cells.Item(46,1).Value = "Kaitoke at Vector Gas Line"
$ws.Cells.Item(46,2).Value = "E coli (Median)"
$ws.Cells.Item(46,3).Value = "D"
$ws.Cells.Item(46,4).Value = "2019 - 2023"
$ws.Cells.Item(46,5).Value = "RepSite"
$ws.Cells.Item(46,6).Value = 250
$ws.Cells.Item(46,7).Value = 931.722340626303
$ws.Cells.Item(46,8).Value = 12981.0063938204
$ws.Cells.Item(46,9).Value = 5860
$ws.Cells.Item(46,10).Value = 27.7777777777778
$ws.Cells.Item(46,11).Value = 48.1481481481481
$ws.Cells.Item(46,12).Value = 222
$ws.Cells.Item(46,13).Value = 740.16
$ws.Cells.Item(46,14).Value = 2696
$ws.Cells.Item(46,15).Value = 1773468
$ws.Cells.Item(46,16).Value = 5573594
$ws.Cells.Item(46,17).Value = "Whanganui District"
$ws.Cells.Item(46,18).Value = "Whanganui"
$ws.Cells.Item(46,19).Value = "Kaitoke Lakes"
$ws.Cells.Item(46,20).Value = "West_4"
$ws.Cells.Item(46,21).Value = "E. coli/100 mL"

# Row 47
$ws.Cells.Item(47,1).Value = "Kaitoke at Vector Gas Line"
$ws.Cells.Item(47,2).Value = "E coli (95th Percentile)"
$ws.Cells.Item(47,3).Value = "E"
$ws.Cells.Item(47,4).Value = "2019 - 2023"
$ws.Cells.Item(47,5).Value = "RepSite"
$ws.Cells.Item(47,6).Value = 250
$ws.Cells.Item(47,7).Value = 931.722340626303
$ws.Cells.Item(47,8).Value = 12981.0063938204
$ws.Cells.Item(47,9).Value = 5860
$ws.Cells.Item(47,10).Value = 27.7777777777778
$ws.Cells.Item(47,11).Value = 48.1481481481481
$ws.Cells.Item(47,12).Value = 222
$ws.Cells.Item(47,13).Value = 740.16
$ws.Cells.Item(47,14).Value = 2696
$ws.Cells.Item(47,15).Value = 1773468
$ws.Cells.Item(47,16).Value = 5573594
$ws.Cells.Item(47,17).Value = "Whanganui District"
$ws.Cells.Item(47,18).Value = "Whanganui"
$ws.Cells.Item(47,19).Value = "Kaitoke Lakes"
$ws.Cells.Item(47,20).Value = "West_4"
$ws.Cells.Item(47,21).Value = "E. coli/100 mL"

# Row 48
$ws.Cells.Item(48,1).Value = "Kaitoke at Vector Gas Line"
$ws.Cells.Item(48,2).Value = "Ammoniacal-N (95th Percentile)"
$ws.Cells.Item(48,3).Value = "B"
$ws.Cells.Item(48,4).Value = "2019 - 2023"
$ws.Cells.Item(48,5).Value = "RepSite"
$ws.Cells.Item(48,6).Value = 0.0175
$ws.Cells.Item(48,7).Value = 0.0361994486818432
$ws.Cells.Item(48,8).Value = 0.5
$ws.Cells.Item(48,9).Value = 0.0821
$ws.Cells.Item(48,12).Value = 0.00943
$ws.Cells.Item(48,13).Value = 0.0622
$ws.Cells.Item(48,14).Value = 0.07471
$ws.Cells.Item(48,15).Value = 1773468
$ws.Cells.Item(48,16).Value = 5573594
$ws.Cells.Item(48,17).Value = "Whanganui District"
$ws.Cells.Item(48,18).Value = "Whanganui"
$ws.Cells.Item(48,19).Value = "Kaitoke Lakes"
$ws.Cells.Item(48,20).Value = "West_4"
$ws.Cells.Item(48,21).Value = "mg NH4-N/L"

# Row 49
$ws.Cells.Item(49,1).Value = "Kaitoke at Vector Gas Line"
$ws.Cells.Item(49,2).Value = "Ammoniacal-N (Median)"
$ws.Cells.Item(49,3).Value = "A"
$ws.Cells.Item(49,4).Value = "2019 - 2023"
$ws.Cells.Item(49,5).Value = "RepSite"
$ws.Cells.Item(49,6).Value = 0.0175
$ws.Cells.Item(49,7).Value = 0.0361994486818432
$ws.Cells.Item(49,8).Value = 0.5
$ws.Cells.Item(49,9).Value = 0.0821
$ws.Cells.Item(49,12).Value = 0.00943
$ws.Cells.Item(49,13).Value = 0.0622
$ws.Cells.Item(49,14).Value = 0.07471
$ws.Cells.Item(49,15).Value = 1773468
$ws.Cells.Item(49,16).Value = 5573594
$ws.Cells.Item(49,17).Value = "Whanganui District"
$ws.Cells.Item(49,18).Value = "Whanganui"
$ws.Cells.Item(49,19).Value = "Kaitoke Lakes"
$ws.Cells.Item(49,20).Value = "West_4"
$ws.Cells.Item(49,21).Value = "mg NH4-N/L"

# Row 50
$ws.Cells.Item(50,1).Value = "Kaitoke at Vector Gas Line"
$ws.Cells.Item(50,2).Value = "Nitrate-N (95th Percentile)"
$ws.Cells.Item(50,3).Value = "B"
$ws.Cells.Item(50,4).Value = "2019 - 2023"
$ws.Cells.Item(50,5).Value = "RepSite"
$ws.Cells.Item(50,6).Value = 0.353
$ws.Cells.Item(50,7).Value = 0.494415203353902
$ws.Cells.Item(50,8).Value = 2.31
$ws.Cells.Item(50,9).Value = 1.704
$ws.Cells.Item(50,12).Value = 0.147
$ws.Cells.Item(50,13).Value = 0.77772
$ws.Cells.Item(50,14).Value = 1.2606
$ws.Cells.Item(50,15).Value = 1773468
$ws.Cells.Item(50,16).Value = 5573594
$ws.Cells.Item(50,17).Value = "Whanganui District"
$ws.Cells.Item(50,18).Value = "Whanganui"
$ws.Cells.Item(50,19).Value = "Kaitoke Lakes"
$ws.Cells.Item(50,20).Value = "West_4"
$ws.Cells.Item(50,21).Value = "mg NO3-N/L"

# Row 51
$ws.Cells.Item(51,1).Value = "Kaitoke at Vector Gas Line"
$ws.Cells.Item(51,2).Value = "Nitrate-N (Median)"
$ws.Cells.Item(51,3).Value = "A"
$ws.Cells.Item(51,4).Value = "2019 - 2023"
$ws.Cells.Item(51,5).Value = "RepSite"
$ws.Cells.Item(51,6).Value = 0.353
$ws.Cells.Item(51,7).Value = 0.494415203353902
$ws.Cells.Item(51,8).Value = 2.31
$ws.Cells.Item(51,9).Value = 1.704
$ws.Cells.Item(51,12).Value = 0.147
$ws.Cells.Item(51,13).Value = 0.77772
$ws.Cells.Item(51,14).Value = 1.2606
$ws.Cells.Item(51,15).Value = 1773468
$ws.Cells.Item(51,16).Value = 5573594
$ws.Cells.Item(51,17).Value = "Whanganui District"
$ws.Cells.Item(51,18).Value = "Whanganui"
$ws.Cells.Item(51,19).Value = "Kaitoke Lakes"
$ws.Cells.Item(51,20).Value = "West_4"
$ws.Cells.Item(51,21).Value = "mg NO3-N/L"

# Row 52
$ws.Cells.Item(52,1).Value = "Kaitoke at Vector Gas Line"
$ws.Cells.Item(52,2).Value = "Soluble Inorganic Nitrogen (95th Percentile)"
$ws.Cells.Item(52,4).Value = "2019 - 2023"
$ws.Cells.Item(52,5).Value = "RepSite"
$ws.Cells.Item(52,6).Value = 0.42
$ws.Cells.Item(52,7).Value = 0.559932782505635
$ws.Cells.Item(52,8).Value = 2.43
$ws.Cells.Item(52,9).Value = 1.88
$ws.Cells.Item(52,12).Value = 0.23
$ws.Cells.Item(52,13).Value = 0.9092
$ws.Cells.Item(52,14).Value = 1.4442
$ws.Cells.Item(52,15).Value = 1773468
$ws.Cells.Item(52,16).Value = 5573594
$ws.Cells.Item(52,17).Value = "Whanganui District"
$ws.Cells.Item(52,18).Value = "Whanganui"
$ws.Cells.Item(52,19).Value = "Kaitoke Lakes"
$ws.Cells.Item(52,20).Value = "West_4"
$ws.Cells.Item(52,21).Value = "g/m3"

# Row 53
$ws.Cells.Item(53,1).Value = "Kaitoke at Vector Gas Line"
$ws.Cells.Item(53,2).Value = "Soluble Inorganic Nitrogen (Median)"
$ws.Cells.Item(53,4).Value = "2019 - 2023"
$ws.Cells.Item(53,5).Value = "RepSite"
$ws.Cells.Item(53,6).Value = 0.42
$ws.Cells.Item(53,7).Value = 0.559932782505635
$ws.Cells.Item(53,8).Value = 2.43
$ws.Cells.Item(53,9).Value = 1.88
$ws.Cells.Item(53,12).Value = 0.23
$ws.Cells.Item(53,13).Value = 0.9092
$ws.Cells.Item(53,14).Value = 1.4442
$ws.Cells.Item(53,15).Value = 1773468
$ws.Cells.Item(53,16).Value = 5573594
$ws.Cells.Item(53,17).Value = "Whanganui District"
$ws.Cells.Item(53,18).Value = "Whanganui"
$ws.Cells.Item(53,19).Value = "Kaitoke Lakes"
$ws.Cells.Item(53,20).Value = "West_4"
$ws.Cells.Item(53,21).Value = "g/m3"

# Row 54
$ws.Cells.Item(54,1).Value = "Kaitoke at Vector Gas Line"
$ws.Cells.Item(54,2).Value = "Total Nitrogen (95th Percentile)"
$ws.Cells.Item(54,4).Value = "2019 - 2023"
$ws.Cells.Item(54,5).Value = "RepSite"
$ws.Cells.Item(54,6).Value = 1.29
$ws.Cells.Item(54,7).Value = 1.46740740740741
$ws.Cells.Item(54,8).Value = 3.6
$ws.Cells.Item(54,9).Value = 2.996
$ws.Cells.Item(54,12).Value = 1
$ws.Cells.Item(54,13).Value = 2.0128
$ws.Cells.Item(54,14).Value = 2.7636
$ws.Cells.Item(54,15).Value = 1773468
$ws.Cells.Item(54,16).Value = 5573594
$ws.Cells.Item(54,17).Value = "Whanganui District"
$ws.Cells.Item(54,18).Value = "Whanganui"
$ws.Cells.Item(54,19).Value = "Kaitoke Lakes"
$ws.Cells.Item(54,20).Value = "West_4"
$ws.Cells.Item(54,21).Value = "g/m3"

# Row 55
$ws.Cells.Item(55,1).Value = "Kaitoke at Vector Gas Line"
$ws.Cells.Item(55,2).Value = "Total Nitrogen (Median)"
$ws.Cells.Item(55,4).Value = "2019 - 2023"
$ws.Cells.Item(55,5).Value = "RepSite"
$ws.Cells.Item(55,6).Value = 1.29
$ws.Cells.Item(55,7).Value = 1.46740740740741
$ws.Cells.Item(55,8).Value = 3.6
$ws.Cells.Item(55,9).Value = 2.996
$ws.Cells.Item(55,12).Value = 1
$ws.Cells.Item(55,13).Value = 2.0128
$ws.Cells.Item(55,14).Value = 2.7636
$ws.Cells.Item(55,15).Value = 1773468
$ws.Cells.Item(55,16).Value = 5573594
$ws.Cells.Item(55,17).Value = "Whanganui District"
$ws.Cells.Item(55,18).Value = "Whanganui"
$ws.Cells.Item(55,19).Value = "Kaitoke Lakes"
$ws.Cells.Item(55,20).Value = "West_4"
$ws.Cells.Item(55,21).Value = "g/m3"

# Row 56
$ws.Cells.Item(56,1).Value = "Kaitoke at Vector Gas Line"
$ws.Cells.Item(56,2).Value = "Total Phosphorus (95th Percentile)"
$ws.Cells.Item(56,4).Value = "2019 - 2023"
$ws.Cells.Item(56,5).Value = "RepSite"
$ws.Cells.Item(56,6).Value = 0.138
$ws.Cells.Item(56,7).Value = 0.158314814814815
$ws.Cells.Item(56,8).Value = 0.523
$ws.Cells.Item(56,9).Value = 0.3444
$ws.Cells.Item(56,12).Value = 0.135
$ws.Cells.Item(56,13).Value = 0.22544
$ws.Cells.Item(56,14).Value = 0.34036
$ws.Cells.Item(56,15).Value = 1773468
$ws.Cells.Item(56,16).Value = 5573594
$ws.Cells.Item(56,17).Value = "Whanganui District"
$ws.Cells.Item(56,18).Value = "Whanganui"
$ws.Cells.Item(56,19).Value = "Kaitoke Lakes"
$ws.Cells.Item(56,20).Value = "West_4"
$ws.Cells.Item(56,21).Value = "g/m3"

# Row 57
$ws.Cells.Item(57,1).Value = "Kaitoke at Vector Gas Line"
$ws.Cells.Item(57,2).Value = "Total Phosphorus (Median)"
$ws.Cells.Item(57,4).Value = "2019 - 2023"
$ws.Cells.Item(57,5).Value = "RepSite"
$ws.Cells.Item(57,6).Value = 0.138
$ws.Cells.Item(57,7).Value = 0.158314814814815
$ws.Cells.Item(57,8).Value = 0.523
$ws.Cells.Item(57,9).Value = 0.3444
$ws.Cells.Item(57,12).Value = 0.135
$ws.Cells.Item(57,13).Value = 0.22544
$ws.Cells.Item(57,14).Value = 0.34036
$ws.Cells.Item(57,15).Value = 1773468
$ws.Cells.Item(57,16).Value = 5573594
$ws.Cells.Item(57,17).Value = "Whanganui District"
$ws.Cells.Item(57,18).Value = "Whanganui"
$ws.Cells.Item(57,19).Value = "Kaitoke Lakes"
$ws.Cells.Item(57,20).Value = "West_4"
$ws.Cells.Item(57,21).Value = "g/m3"
